$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = [double]"-4.30813814202642e-9"
$ws.Range("C5").Value = [double]"-4.3081376419921005e-9"
$ws.Range("D5").Value = [double]"-4.308137523078548e-9"
$ws.Range("E5").Value = [double]"-4.3081371509771066e-9"
$ws.Range("F5").Value = [double]"-9.090908574732907e-9"
$ws.Range("B6").Value = [double]"0.04267407439113609"
$ws.Range("C6").Value = [double]"0.042894280744237914"
$ws.Range("D6").Value = [double]"0.0430276985218409"
$ws.Range("E6").Value = [double]"0.043101077761150575"
$ws.Range("F6").Value = [double]"0.042260660261863925"
$ws.Range("B7").Value = [double]"-0.05333160688517414"
$ws.Range("C7").Value = [double]"-0.053288804082253866"
$ws.Range("D7").Value = [double]"-0.05557683628089839"
$ws.Range("E7").Value = [double]"-0.06525313270591816"
$ws.Range("F7").Value = [double]"-0.030267369093617522"
$ws.Range("B8").Value = [double]"0.010657528185899915"
$ws.Range("C8").Value = [double]"0.010394519029878308"
$ws.Range("D8").Value = [double]"0.012549133450919972"
$ws.Range("E8").Value = [double]"0.022152050636630435"
$ws.Range("F8").Value = [double]"-0.011993300259154982"
$ws.Range("B9").Value = [double]"0.005000571413257865"
$ws.Range("C9").Value = [double]"0.0051806039760923725"
$ws.Range("D9").Value = [double]"0.0052472829146281595"
$ws.Range("E9").Value = [double]"0.005307294056864217"
$ws.Range("F9").Value = [double]"0.0052139434288796084"
$ws.Range("B10").Value = [double]"0.009998580069983257"
$ws.Range("C10").Value = [double]"0.01035856792971497"
$ws.Range("D10").Value = [double]"0.010491901005877567"
$ws.Range("E10").Value = [double]"0.010611904026549646"
$ws.Range("F10").Value = [double]"0.010425251031502618"
$ws.Range("B11").Value = [double]"0.017670610425920653"
$ws.Range("C11").Value = [double]"0.016990758226373782"
$ws.Range("D11").Value = [double]"0.01679080757647534"
$ws.Range("E11").Value = [double]"0.016564167637040973"
$ws.Range("F11").Value = [double]"0.016190633876931777"
$ws.Range("B12").Value = [double]"0.005000668208106999"
$ws.Range("C12").Value = [double]"0.005180707863493138"
$ws.Range("D12").Value = [double]"0.005247389492650469"
$ws.Range("E12").Value = [double]"0.005307403085157084"
$ws.Range("F12").Value = [double]"0.0052140486499064305"
$ws.Range("B13").Value = [double]"0.005000630883271751"
$ws.Range("C13").Value = [double]"0.005180667802562927"
$ws.Range("D13").Value = [double]"0.005247348393776203"
$ws.Range("E13").Value = [double]"0.005307361040862082"
$ws.Range("F13").Value = [double]"0.005214008072518304"
$ws.Range("B14").Value = [double]"0.0050003093053147965"
$ws.Range("C14").Value = [double]"0.005180322652479597"
$ws.Range("D14").Value = [double]"0.005246994301581261"
$ws.Range("E14").Value = [double]"0.005306998803677855"
$ws.Range("F14").Value = [double]"0.005213658473036407"
$ws.Range("B15").Value = [double]"0.02397957211052107"
$ws.Range("C15").Value = [double]"0.024520262607329243"
$ws.Range("D15").Value = [double]"0.024690511467573107"
$ws.Range("E15").Value = [double]"0.024811632614961734"
$ws.Range("F15").Value = [double]"0.02441041081199054"
$ws.Range("B16").Value = [double]"0.18514854580690884"
$ws.Range("C16").Value = [double]"0.18528930995297396"
$ws.Range("D16").Value = [double]"0.18419856393761458"
$ws.Range("E16").Value = [double]"0.1801950824656111"
$ws.Range("F16").Value = [double]"-0.20338184217921984"
$ws.Range("B17").Value = [double]"0.19489585424575645"
$ws.Range("C17").Value = [double]"0.19503315331446322"
$ws.Range("D17").Value = [double]"0.19466177812698635"
$ws.Range("E17").Value = [double]"0.19345732332778878"
$ws.Range("F17").Value = [double]"-0.20634294306688447"
$ws.Range("B18").Value = [double]"0.0025000760121757044"
$ws.Range("C18").Value = [double]"0.0025900909121572874"
$ws.Range("D18").Value = [double]"0.002623429978589453"
$ws.Range("E18").Value = [double]"0.0026534352373280925"
$ws.Range("F18").Value = [double]"0.0026067604286568995"
$ws.Range("B19").Value = [double]"0.005485414961530933"
$ws.Range("C19").Value = [double]"0.0056519068619503"
$ws.Range("D19").Value = [double]"0.005713549235630298"
$ws.Range("E19").Value = [double]"0.005769239685276451"
$ws.Range("F19").Value = [double]"0.005684609531719483"
$ws.Range("B20").Value = [double]"0.01097989295530378"
$ws.Range("C20").Value = [double]"0.011083805500507771"
$ws.Range("D20").Value = [double]"0.011092287501341285"
$ws.Range("E20").Value = [double]"0.011067638298749816"
$ws.Range("F20").Value = [double]"0.010891591649661435"
$ws.Range("B21").Value = [double]"0.0025001738319323993"
$ws.Range("C21").Value = [double]"0.002590195899554782"
$ws.Range("D21").Value = [double]"0.002623537685094424"
$ws.Range("E21").Value = [double]"0.002653545420043059"
$ws.Range("F21").Value = [double]"0.0026068667637804775"
$ws.Range("B22").Value = [double]"0.0025001978573772525"
$ws.Range("C22").Value = [double]"0.002590221685966609"
$ws.Range("D22").Value = [double]"0.0026235641395274552"
$ws.Range("E22").Value = [double]"0.002653572482946556"
$ws.Range("F22").Value = [double]"0.0026068928825852804"
$ws.Range("B23").Value = [double]"0.002499871678145875"
$ws.Range("C23").Value = [double]"0.0025898715973288276"
$ws.Range("D23").Value = [double]"0.002623204980830519"
$ws.Range("E23").Value = [double]"0.0026532050627183566"
$ws.Range("F23").Value = [double]"0.002606538280885042"
$ws.Range("B24").Value = [double]"0.002258550464868548"
$ws.Range("C24").Value = [double]"0.0023010295875199606"
$ws.Range("D24").Value = [double]"0.0023197321929759756"
$ws.Range("E24").Value = [double]"0.0023313410357977448"
$ws.Range("F24").Value = [double]"0.002245104412455492"
$ws.Range("B25").Value = [double]"0.034993096112758186"
$ws.Range("C25").Value = [double]"0.03503794560576232"
$ws.Range("D25").Value = [double]"0.03489285229035758"
$ws.Range("E25").Value = [double]"0.03461982326703022"
$ws.Range("F25").Value = [double]"0.03985314874887429"
$ws.Range("B26").Value = [double]"0.03591092899623125"
$ws.Range("C26").Value = [double]"0.03595617374290462"
$ws.Range("D26").Value = [double]"0.03586635488843895"
$ws.Range("E26").Value = [double]"0.035739820001720116"
$ws.Range("F26").Value = [double]"0.04026879686539362"
$ws.Range("B27").Value = [double]"2.9477240337494127e-5"
$ws.Range("C27").Value = [double]"3.1638219447170265e-5"
$ws.Range("D27").Value = [double]"3.245798453529442e-5"
$ws.Range("E27").Value = [double]"3.3204719908694875e-5"
$ws.Range("F27").Value = [double]"3.204659052801731e-5"
$ws.Range("B28").Value = [double]"0.00012265919940622017"
$ws.Range("C28").Value = [double]"0.00013132019027135063"
$ws.Range("D28").Value = [double]"0.00013460314228443745"
$ws.Range("E28").Value = [double]"0.00013759490122489384"
$ws.Range("F28").Value = [double]"0.00013297623962034557"
$ws.Range("B29").Value = [double]"0.00040817636751258406"
$ws.Range("C29").Value = [double]"0.00038811801583059195"
$ws.Range("D29").Value = [double]"0.0003819261832578754"
$ws.Range("E29").Value = [double]"0.0003742823131069778"
$ws.Range("F29").Value = [double]"0.00035909421012304837"
$ws.Range("B30").Value = [double]"2.94925839774373e-5"
$ws.Range("C30").Value = [double]"3.165444199722731e-5"
$ws.Range("D30").Value = [double]"3.247454541806554e-5"
$ws.Range("E30").Value = [double]"3.32215428713223e-5"
$ws.Range("F30").Value = [double]"3.206250744860781e-5"
$ws.Range("B31").Value = [double]"2.9492345111550528e-5"
$ws.Range("C31").Value = [double]"3.16541763908656e-5"
$ws.Range("D31").Value = [double]"3.24742694227237e-5"
$ws.Range("E31").Value = [double]"3.3221257298421744e-5"
$ws.Range("F31").Value = [double]"3.206223669237753e-5"
$ws.Range("B32").Value = [double]"2.9492345111550528e-5"
$ws.Range("C32").Value = [double]"3.16541763908656e-5"
$ws.Range("D32").Value = [double]"3.24742694227237e-5"
$ws.Range("E32").Value = [double]"3.322125729842175e-5"
$ws.Range("F32").Value = [double]"3.206223669237752e-5"
$ws.Range("B34").Value = [double]"1.060346700297827"
$ws.Range("C34").Value = [double]"1.0603387003608413"
$ws.Range("D34").Value = [double]"1.060335795073358"
$ws.Range("E34").Value = [double]"1.0603338580037491"
$ws.Range("F34").Value = [double]"1.0603438759658415"
$ws.Range("B35").Value = [double]"1.0582828207010737"
$ws.Range("C35").Value = [double]"1.0582801620312687"
$ws.Range("D35").Value = [double]"1.058316516381771"
$ws.Range("E35").Value = [double]"1.0584602553159828"
$ws.Range("F35").Value = [double]"1.0644961639565997"
$ws.Range("B36").Value = [double]"1.0601823565297992"
$ws.Range("C36").Value = [double]"1.0601684399518805"
$ws.Range("D36").Value = [double]"1.0601633433112234"
$ws.Range("E36").Value = [double]"1.0601594340216"
$ws.Range("F36").Value = [double]"1.0601725198806582"
$ws.Range("B37").Value = [double]"1.0602186842782253"
$ws.Range("C37").Value = [double]"1.0602064581406496"
$ws.Range("D37").Value = [double]"1.0602019878344746"
$ws.Range("E37").Value = [double]"1.0601986396096783"
$ws.Range("F37").Value = [double]"1.0602108279146962"
$ws.Range("B38").Value = [double]"1.0598444592499354"
$ws.Range("C38").Value = [double]"1.0598464748760148"
$ws.Range("D38").Value = [double]"1.0598469056179907"
$ws.Range("E38").Value = [double]"1.0598493556107205"
$ws.Range("F38").Value = [double]"1.0598689981487588"
$ws.Range("B39").Value = [double]"1.0596307922326127"
$ws.Range("C39").Value = [double]"1.0596251153976626"
$ws.Range("D39").Value = [double]"1.0596226970720841"
$ws.Range("E39").Value = [double]"1.059622582900386"
$ws.Range("F39").Value = [double]"1.059646214137439"
$ws.Range("B40").Value = [double]"1.05968011149349"
$ws.Range("C40").Value = [double]"1.0596762101863413"
$ws.Range("D40").Value = [double]"1.0596744494642596"
$ws.Range("E40").Value = [double]"1.0596749271359869"
$ws.Range("F40").Value = [double]"1.0596976377278038"
$ws.Range("B41").Value = [double]"1.0596393744476884"
$ws.Range("C41").Value = [double]"1.0596340065837166"
$ws.Range("D41").Value = [double]"1.059631702692032"
$ws.Range("E41").Value = [double]"1.0596316915109292"
$ws.Range("F41").Value = [double]"1.0596551625404196"
$ws.Range("B42").Value = [double]"1.057720892318825"
$ws.Range("C42").Value = [double]"1.0577227911701494"
$ws.Range("D42").Value = [double]"1.0576955928437255"
$ws.Range("E42").Value = [double]"1.0575674196332712"
$ws.Range("F42").Value = [double]"1.0642859594585514"
$ws.Range("B43").Value = [double]"-0.0008908827352805736"
$ws.Range("C43").Value = [double]"-0.0006311047860027225"
$ws.Range("D43").Value = [double]"-0.0005052228634117533"
$ws.Range("E43").Value = [double]"-0.00035943700716228546"
$ws.Range("F43").Value = [double]"-0.00036800863153772095"
$ws.Range("B44").Value = [double]"1.438684068520384e-19"
$ws.Range("C44").Value = [double]"-1.0963699529388172e-19"
$ws.Range("D44").Value = [double]"1.497565535865285e-19"
$ws.Range("E44").Value = [double]"8.922480328977489e-20"
$ws.Range("F44").Value = [double]"8.922480328977489e-20"
$ws.Range("B45").Value = [double]"-0.00048440799540957796"
$ws.Range("C45").Value = [double]"-0.0004708374462980982"
$ws.Range("D45").Value = [double]"-0.00046578955991228996"
$ws.Range("E45").Value = [double]"-0.000461458546531027"
$ws.Range("F45").Value = [double]"-0.00047019489609798017"
$ws.Range("B46").Value = [double]"-8.477988412594885e-5"
$ws.Range("C46").Value = [double]"-8.859281628432011e-5"
$ws.Range("D46").Value = [double]"-8.965724198303109e-5"
$ws.Range("E46").Value = [double]"-9.084239596753449e-5"
$ws.Range("F46").Value = [double]"-9.305465563245771e-5"
$ws.Range("B47").Value = [double]"7.877236809590782e-6"
$ws.Range("C47").Value = [double]"7.903491312090182e-6"
$ws.Range("D47").Value = [double]"7.916606733211724e-6"
$ws.Range("E47").Value = [double]"7.929570743402394e-6"
$ws.Range("F47").Value = [double]"7.920311578617106e-6"
$ws.Range("B48").Value = [double]"9.900110659230514e-6"
$ws.Range("C48").Value = [double]"9.892543046976181e-6"
$ws.Range("D48").Value = [double]"9.887774846979489e-6"
$ws.Range("E48").Value = [double]"9.877221582804596e-6"
$ws.Range("F48").Value = [double]"9.8369666787723e-6"
$ws.Range("B49").Value = [double]"7.167144649691226e-6"
$ws.Range("C49").Value = [double]"7.141663629376297e-6"
$ws.Range("D49").Value = [double]"7.131605279617561e-6"
$ws.Range("E49").Value = [double]"7.120572681687092e-6"
$ws.Range("F49").Value = [double]"7.130377261789335e-6"
$ws.Range("B50").Value = [double]"8.897892261574474e-6"
$ws.Range("C50").Value = [double]"8.904198979781336e-6"
$ws.Range("D50").Value = [double]"8.908311026319314e-6"
$ws.Range("E50").Value = [double]"8.917285645474704e-6"
$ws.Range("F50").Value = [double]"8.950953762173813e-6"
$ws.Range("B51").Value = [double]"0.011663689012179963"
$ws.Range("C51").Value = [double]"0.011663679788263947"
$ws.Range("D51").Value = [double]"0.011663693611734584"
$ws.Range("E51").Value = [double]"0.011663731364275882"
$ws.Range("B52").Value = [double]"0.005833663902745664"
$ws.Range("C52").Value = [double]"0.005833688977608943"
$ws.Range("D52").Value = [double]"0.005833705194212184"
$ws.Range("E52").Value = [double]"0.005833701938246927"
